# Generate Report for Handback
# - The "20a5c65a-f8fb-499a-a53a-41282ae29a11" file's status flips from
#   "Ready for handoff" to "Handback transform failed" (Overview + both
#   language sheets), and each language sheet gets a detailed error message
#   in its "Error Detail" column (P) for that row, explaining the handback
#   file-name mismatch. Column P is also widened to fit the longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: zh-cn / de-de status columns (E3, F3) for the
# 20a5c65a-f8fb-499a-a53a-41282ae29a11.md row
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# zh-cn / de-de detail sheets: Status column (C3) for the same row
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# zh-cn / de-de detail sheets: Error Detail column (P3) gets the transform
# failure explanation
$wsZhCn.Range("P3").Value = "Handback file name: ddsc0upu.nf1 is different with handoff file name: 20a5c65a-f8fb-499a-a53a-41282ae29a11.e8ff649b42034cd2db58cd33a71d9cb98689a775.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: ddsc0upu.nf1 is different with handoff file name: 20a5c65a-f8fb-499a-a53a-41282ae29a11.e8ff649b42034cd2db58cd33a71d9cb98689a775.de-de."

# Widen the "Error Detail" column (P) on both sheets so the new message is
# readable (stored OOXML column width of 40 characters).
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
